# Automatic update of files.
# Increment the "Förändrad" (Changed) date in column C (rows 2-43) by one day,
# i.e. from serial 45832 (2025-06-24) to 45833 (2025-06-25).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 43 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -ne $null) {
        $cell.Value = $cell.Value2 + 1
    }
}
